$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 737.25
$ws.Range("I12").Value = 649.5
$ws.Range("K12").Value = 649.5
$ws.Range("M12").Value = -479.5
$ws.Range("H17").Value = 2211.5715
$ws.Range("J17").Value = 2272.15
$ws.Range("L17").Value = 6816.450000000001
$ws.Range("N17").Value = -7152.450000000001
$ws.Range("H40").Value = 4685.6665
$ws.Range("J40").Value = 7406.6665
$ws.Range("L40").Value = 7406.6665
$ws.Range("N40").Value = -7756.6665
$ws.Range("H94").Value = 3134.2856
$ws.Range("I94").Value = 3134.2856
$ws.Range("K94").Value = 3134.2856
$ws.Range("M94").Value = -2683.2856
$ws.Range("H98").Value = 1568.6666
$ws.Range("I98").Value = 1666.6666
$ws.Range("J98").Value = 1470.6666
$ws.Range("K98").Value = 1666.6666
$ws.Range("L98").Value = 1470.6666
$ws.Range("M98").Value = -168.6666
$ws.Range("N98").Value = -4466.6666
$ws.Range("H112").Value = 3276
$ws.Range("J112").Value = 3488.6667
$ws.Range("L112").Value = 10466.0001
$ws.Range("N112").Value = -12682.0001
$ws.Range("H116").Value = 4283
$ws.Range("I116").Value = 4175
$ws.Range("J116").Value = 4499
$ws.Range("K116").Value = 4175
$ws.Range("L116").Value = 4499
$ws.Range("M116").Value = -733
$ws.Range("N116").Value = -11383
$ws.Range("H122").Value = 1568.6666
$ws.Range("I122").Value = 1666.6666
$ws.Range("J122").Value = 1470.6666
$ws.Range("K122").Value = 4999.9998
$ws.Range("L122").Value = 4411.9998
$ws.Range("M122").Value = -2549.9998
$ws.Range("N122").Value = -9311.9998
$ws.Range("H125").Value = 1677
$ws.Range("I125").Value = 816
$ws.Range("J125").Value = 2538
$ws.Range("K125").Value = 7344
$ws.Range("L125").Value = 22842
$ws.Range("M125").Value = -4884
$ws.Range("N125").Value = -27762
$ws.Range("H132").Value = 14515.866
$ws.Range("I132").Value = 14515.866
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 43547.598
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -41017.598
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 1777.091

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7658.3335
$ws.Range("I74").Value = 7386.952
$ws.Range("J74").Value = 8608.166999999999
$ws.Range("K74").Value = 7386.952
$ws.Range("L74").Value = 8608.166999999999
$ws.Range("M74").Value = -6512.952
$ws.Range("N74").Value = -10356.167
$ws.Range("H77").Value = 7658.3335
$ws.Range("I77").Value = 7386.952
$ws.Range("J77").Value = 8608.166999999999
$ws.Range("K77").Value = 36934.76
$ws.Range("L77").Value = 43040.835
$ws.Range("M77").Value = -32566.76
$ws.Range("N77").Value = -51776.835
$ws.Range("H88").Value = 1072.2727
$ws.Range("I88").Value = 447.5
$ws.Range("J88").Value = 1429.2858
$ws.Range("K88").Value = 447.5
$ws.Range("L88").Value = 1429.2858
$ws.Range("M88").Value = -41.5
$ws.Range("N88").Value = -2241.2858
$ws.Range("H91").Value = 1072.2727
$ws.Range("I91").Value = 447.5
$ws.Range("J91").Value = 1429.2858
$ws.Range("K91").Value = 447.5
$ws.Range("L91").Value = 1429.2858
$ws.Range("M91").Value = 956.5
$ws.Range("N91").Value = -4237.2858
$ws.Range("H97").Value = 1570
$ws.Range("I97").Value = 1246.5
$ws.Range("J97").Value = 2055.25
$ws.Range("K97").Value = 1246.5
$ws.Range("L97").Value = 2055.25
$ws.Range("M97").Value = -750.5
$ws.Range("N97").Value = -3047.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2480.5151
$ws.Range("I86").Value = 1440.3846
$ws.Range("J86").Value = 6343.857
$ws.Range("K86").Value = 1440.3846
$ws.Range("L86").Value = 6343.857
$ws.Range("M86").Value = -317.3846000000001
$ws.Range("N86").Value = -8589.857
$ws.Range("H89").Value = 2480.5151
$ws.Range("I89").Value = 1440.3846
$ws.Range("J89").Value = 6343.857
$ws.Range("K89").Value = 7201.923000000001
$ws.Range("L89").Value = 31719.285
$ws.Range("M89").Value = -1585.923000000001
$ws.Range("N89").Value = -42951.285
$ws.Range("H94").Value = 773.75
$ws.Range("I94").Value = 747.5
$ws.Range("K94").Value = 747.5
$ws.Range("M94").Value = -296.5
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H105").Value = 1917.6
$ws.Range("I105").Value = 1917.6
$ws.Range("K105").Value = 1917.6
$ws.Range("M105").Value = -170.5999999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 456.57144
$ws.Range("I6").Value = 439.4
$ws.Range("K6").Value = 439.4
$ws.Range("M6").Value = -326.4
$ws.Range("H31").Value = 5495.405
$ws.Range("I31").Value = 2834.625
$ws.Range("J31").Value = 7132.8076
$ws.Range("K31").Value = 2834.625
$ws.Range("L31").Value = 7132.8076
$ws.Range("M31").Value = -2539.625
$ws.Range("N31").Value = -7722.8076
$ws.Range("H34").Value = 5495.405
$ws.Range("I34").Value = 2834.625
$ws.Range("J34").Value = 7132.8076
$ws.Range("K34").Value = 2834.625
$ws.Range("L34").Value = 7132.8076
$ws.Range("M34").Value = -2632.625
$ws.Range("N34").Value = -7536.8076

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 194.78572
$ws.Range("J12").Value = 253.9
$ws.Range("L12").Value = 761.7
$ws.Range("N12").Value = -1107.7
$ws.Range("H39").Value = 8156
$ws.Range("J39").Value = 9206.857
$ws.Range("L39").Value = 27620.571
$ws.Range("N39").Value = -28208.571
$ws.Range("H55").Value = 5769
$ws.Range("J55").Value = 6612.4
$ws.Range("L55").Value = 19837.2
$ws.Range("N55").Value = -20191.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 889.5
$ws.Range("I97").Value = 875.1429000000001
$ws.Range("K97").Value = 875.1429000000001
$ws.Range("M97").Value = -379.1429000000001
$ws.Range("H122").Value = 3289.4
$ws.Range("I122").Value = 2662.75
$ws.Range("J122").Value = 5796
$ws.Range("K122").Value = 7988.25
$ws.Range("L122").Value = 17388
$ws.Range("M122").Value = -5538.25
$ws.Range("N122").Value = -22288
$ws.Range("H134").Value = 97966.71000000001
$ws.Range("J134").Value = 97966.71000000001
$ws.Range("L134").Value = 293900.13
$ws.Range("N134").Value = -298970.13

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4173.2666
$ws.Range("I61").Value = 2289
$ws.Range("K61").Value = 2289
$ws.Range("M61").Value = -2087
$ws.Range("H87").Value = 30000
$ws.Range("I87").Value = 30000
$ws.Range("K87").Value = 30000
$ws.Range("M87").Value = -28877
$ws.Range("H90").Value = 30000
$ws.Range("I90").Value = 30000
$ws.Range("K90").Value = 90000
$ws.Range("M90").Value = -84384
$ws.Range("H101").Value = 22465
$ws.Range("J101").Value = 22465
$ws.Range("L101").Value = 22465
$ws.Range("N101").Value = -28955
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 4173.2666
$ws.Range("I113").Value = 2289
$ws.Range("K113").Value = 2289
$ws.Range("M113").Value = -119
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 633
$ws.Range("I81").Value = 633
$ws.Range("K81").Value = 1266
$ws.Range("M81").Value = -205
$ws.Range("H84").Value = 633
$ws.Range("I84").Value = 633
$ws.Range("K84").Value = 6330
$ws.Range("M84").Value = -1026
$ws.Range("H96").Value = 1223.3334
$ws.Range("I96").Value = 1206.375
$ws.Range("J96").Value = 1242.7142
$ws.Range("K96").Value = 1206.375
$ws.Range("L96").Value = 1242.7142
$ws.Range("M96").Value = 166.625
$ws.Range("N96").Value = -3988.7142
$ws.Range("H136").Value = 1550
$ws.Range("I136").Value = 1100
$ws.Range("K136").Value = 3300
$ws.Range("M136").Value = -750
